# Generate Report for Handback
# Refresh the Handoff/Handback timestamps for 85360c91-fe2d-40fc-b68f-10e2238ac63f.md
# after a new HO Xliff generation / handback report run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-11-09 17:57:25"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-11-09 17:57:12"
$zhcn.Range("K3").Value = "2016-11-09 17:58:04"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-11-09 17:57:25"
$dede.Range("K3").Value = "2016-11-09 17:58:21"
